$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

Set-TextCell "D2" "296.66"
Set-TextCell "E2" "2.98%"
Set-TextCell "G2" "19"
Set-TextCell "D3" "41.18"
Set-TextCell "E3" "2.16%"
Set-TextCell "G3" "19"
Set-TextCell "D4" "5.015"
Set-TextCell "E4" "-0.48%"
Set-TextCell "G4" "19"
Set-TextCell "D5" "0.07500"
Set-TextCell "E5" "2.67%"
Set-TextCell "G5" "19"
Set-TextCell "B6" "FTXToken"
Set-TextCell "C6" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell "D6" "1.572"
Set-TextCell "E6" "4.01%"
Set-TextCell "G6" "19"
Set-TextCell "B7" "MXToken"
Set-TextCell "C7" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D7" "0.9280"
Set-TextCell "E7" "0.86%"
Set-TextCell "G7" "19"
Set-TextCell "B8" "BTSEToken"
Set-TextCell "C8" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell "D8" "2.405"
Set-TextCell "E8" "0.34%"
Set-TextCell "G8" "19"
Set-TextCell "B9" "LiechtensteinCryptoassetsExchange"
Set-TextCell "C9" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D9" "0.1218"
Set-TextCell "E9" "2.50%"
Set-TextCell "G9" "19"
Set-TextCell "B10" "WazirX"
Set-TextCell "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D10" "0.1845"
Set-TextCell "E10" "6.66%"
Set-TextCell "G10" "19"
Set-TextCell "B11" "MandalaExchangeToken"
Set-TextCell "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D11" "0.08905"
Set-TextCell "E11" "3.12%"
Set-TextCell "G11" "19"
Set-TextCell "B12" "BitrueCoin"
Set-TextCell "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D12" "0.04091"
Set-TextCell "E12" "-1.88%"
Set-TextCell "G12" "19"
Set-TextCell "B13" "BitMartToken"
Set-TextCell "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D13" "0.1054"
Set-TextCell "E13" "0.08%"
Set-TextCell "G13" "19"
Set-TextCell "B14" "BitForexToken"
Set-TextCell "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D14" "0.001286"
Set-TextCell "E14" "1.69%"
Set-TextCell "G14" "19"
Set-TextCell "D15" "0.005854"
Set-TextCell "E15" "-0.99%"
Set-TextCell "G15" "19"
Set-TextCell "D16" "3.343"
Set-TextCell "E16" "-1.54%"
Set-TextCell "G16" "19"
Set-TextCell "B17" "GateToken"
Set-TextCell "C17" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell "D17" "4.367"
Set-TextCell "E17" "1.95%"
Set-TextCell "G17" "19"
Set-TextCell "B18" "BitpandaEcosystemToken"
Set-TextCell "C18" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextCell "D18" "0.3287"
Set-TextCell "E18" "-0.04%"
Set-TextCell "G18" "19"
Set-TextCell "B19" "MCDex"
Set-TextCell "C19" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell "D19" "8.002"
Set-TextCell "E19" "2.24%"
Set-TextCell "G19" "19"
Set-TextCell "B20" "ProBitToken"
Set-TextCell "C20" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextCell "D20" "0.1419"
Set-TextCell "E20" "5.13%"
Set-TextCell "G20" "19"
Set-TextCell "B21" "ZBToken"
Set-TextCell "C21" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextCell "D21" "0.2966"
Set-TextCell "E21" "2.97%"
Set-TextCell "G21" "19"
Set-TextCell "D22" "0.04054"
Set-TextCell "E22" "5.03%"
Set-TextCell "G22" "19"
Set-TextCell "E23" "-0.11%"
Set-TextCell "G23" "19"
Set-TextCell "D24" "0.003894"
Set-TextCell "E24" "1.39%"
Set-TextCell "G24" "19"
Set-TextCell "D25" "0.0001231"
Set-TextCell "E25" "-3.87%"
Set-TextCell "G25" "19"
Set-TextCell "E26" "0.07%"
Set-TextCell "G26" "19"
Set-TextCell "G27" "19"
Set-TextCell "G28" "19"
Set-TextCell "G29" "19"
Set-TextCell "G30" "19"
Set-TextCell "G31" "19"
Set-TextCell "G32" "19"
Set-TextCell "G33" "19"
Set-TextCell "G34" "19"
Set-TextCell "G35" "19"
Set-TextCell "G36" "19"
Set-TextCell "G37" "19"
Set-TextCell "D38" "0.02416"
Set-TextCell "E38" "4.33%"
Set-TextCell "G38" "19"
Set-TextCell "D39" "0.05208"
Set-TextCell "E39" "4.56%"
Set-TextCell "G39" "19"
Set-TextCell "D40" "0.005996"
Set-TextCell "E40" "-8.51%"
Set-TextCell "G40" "19"
Set-TextCell "D41" "0.007821"
Set-TextCell "E41" "1.64%"
Set-TextCell "G41" "19"
Set-TextCell "D42" "0.1326"
Set-TextCell "E42" "4.03%"
Set-TextCell "G42" "19"
Set-TextCell "D43" "0.007377"
Set-TextCell "E43" "0.35%"
Set-TextCell "G43" "19"
Set-TextCell "D44" "0.008135"
Set-TextCell "E44" "15.18%"
Set-TextCell "G44" "19"
Set-TextCell "D45" "0.2975"
Set-TextCell "E45" "-4.79%"
Set-TextCell "G45" "19"
Set-TextCell "D46" "0.00006246"
Set-TextCell "E46" "-2.72%"
Set-TextCell "G46" "19"
Set-TextCell "E47" "0.03%"
Set-TextCell "G47" "19"
Set-TextCell "D48" "0.04618"
Set-TextCell "E48" "-81.57%"
Set-TextCell "G48" "19"
Set-TextCell "D49" "0.004200"
Set-TextCell "E49" "-0.03%"
Set-TextCell "G49" "19"
Set-TextCell "E50" "0.03%"
Set-TextCell "G50" "19"
Set-TextCell "E51" "0.03%"
Set-TextCell "G51" "19"
